$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.256564333333333
$ws.Range("H2").Value = 15.769693
$ws.Range("I2").Value = 0.003747859920520347
$ws.Range("J2").Value = 0.003747859920520347
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.970048
$ws.Range("N2").Value = 26.910144
$ws.Range("O2").Value = 0.487108783009476
$ws.Range("P2").Value = 0.4871087830094759
$ws.Range("Q2").Value = 47.151634385088
$ws.Range("R2").Value = 424.364709465792
$ws.Range("S2").Value = 0.001825615484774658
$ws.Range("T2").Value = 0.001825615484774657
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.256564333333333
$ws.Range("H3").Value = 15.769693
$ws.Range("I3").Value = 0.003747859920520347
$ws.Range("J3").Value = 0.003747859920520347
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.012070666666666
$ws.Range("N3").Value = 27.036212
$ws.Range("O3").Value = 0.489390778604016
$ws.Range("P3").Value = 0.489390778604016
$ws.Range("Q3").Value = 47.37252923587955
$ws.Range("R3").Value = 426.352763122916
$ws.Range("S3").Value = 0.001834168084602238
$ws.Range("T3").Value = 0.001834168084602238
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.256564333333333
$ws.Range("H4").Value = 15.769693
$ws.Range("I4").Value = 0.003747859920520347
$ws.Range("J4").Value = 0.003747859920520347
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4327576666666667
$ws.Range("N4").Value = 1.298273
$ws.Range("O4").Value = 0.02350043838650813
$ws.Range("P4").Value = 0.02350043838650813
$ws.Range("Q4").Value = 2.274818515576555
$ws.Range("R4").Value = 20.473366640189
$ws.Range("S4").Value = 0.000088076351143451665466968525
$ws.Range("T4").Value = 0.000088076351143451651914441369
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1312.703450666667
$ws.Range("H5").Value = 3938.110352
$ws.Range("I5").Value = 0.93593996730609
$ws.Range("J5").Value = 0.9359399673060897
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.970048
$ws.Range("N5").Value = 26.910144
$ws.Range("O5").Value = 0.487108783009476
$ws.Range("P5").Value = 0.4871087830094759
$ws.Range("Q5").Value = 11775.01296224563
$ws.Range("R5").Value = 105975.1166602107
$ws.Range("S5").Value = 0.4559045784443982
$ws.Range("T5").Value = 0.455904578444398
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1312.703450666667
$ws.Range("H6").Value = 3938.110352
$ws.Range("I6").Value = 0.93593996730609
$ws.Range("J6").Value = 0.9359399673060897
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.012070666666666
$ws.Range("N6").Value = 27.036212
$ws.Range("O6").Value = 0.489390778604016
$ws.Range("P6").Value = 0.489390778604016
$ws.Range("Q6").Value = 11830.17626178518
$ws.Range("R6").Value = 106471.5863560666
$ws.Range("S6").Value = 0.4580403893265446
$ws.Range("T6").Value = 0.4580403893265445
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1312.703450666667
$ws.Range("H7").Value = 3938.110352
$ws.Range("I7").Value = 0.93593996730609
$ws.Range("J7").Value = 0.9359399673060897
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.4327576666666667
$ws.Range("N7").Value = 1.298273
$ws.Range("O7").Value = 0.02350043838650813
$ws.Range("P7").Value = 0.02350043838650813
$ws.Range("Q7").Value = 568.0824823357884
$ws.Range("R7").Value = 5112.742341022095
$ws.Range("S7").Value = 0.0219949995351472
$ws.Range("T7").Value = 0.02199499953514719
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 84.59089266666666
$ws.Range("H8").Value = 253.772678
$ws.Range("I8").Value = 0.06031217277338979
$ws.Range("J8").Value = 0.06031217277338978
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.970048
$ws.Range("N8").Value = 26.910144
$ws.Range("O8").Value = 0.487108783009476
$ws.Range("P8").Value = 0.4871087830094759
$ws.Range("Q8").Value = 758.7843675828479
$ws.Range("R8").Value = 6829.059308245632
$ws.Range("S8").Value = 0.02937858908030315
$ws.Range("T8").Value = 0.02937858908030314
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 84.59089266666666
$ws.Range("H9").Value = 253.772678
$ws.Range("I9").Value = 0.06031217277338979
$ws.Range("J9").Value = 0.06031217277338978
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.012070666666666
$ws.Range("N9").Value = 27.036212
$ws.Range("O9").Value = 0.489390778604016
$ws.Range("P9").Value = 0.489390778604016
$ws.Range("Q9").Value = 762.339102468415
$ws.Range("R9").Value = 6861.051922215735
$ws.Range("S9").Value = 0.02951622119286916
$ws.Range("T9").Value = 0.02951622119286916
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 84.59089266666666
$ws.Range("H10").Value = 253.772678
$ws.Range("I10").Value = 0.06031217277338979
$ws.Range("J10").Value = 0.06031217277338978
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.4327576666666667
$ws.Range("N10").Value = 1.298273
$ws.Range("O10").Value = 0.02350043838650813
$ws.Range("P10").Value = 0.02350043838650813
$ws.Range("Q10").Value = 36.60735733167711
$ws.Range("R10").Value = 329.466215985094
$ws.Range("S10").Value = 0.00141736250021748
$ws.Range("T10").Value = 0.00141736250021748
